# [V1.4.19] Diller verisi guncellendi.
# The "Diller" (languages) sheet had two stray test/debug rows
# ("history" and "check") accidentally mixed into the otherwise
# alphabetically-sorted list of language names in column A.
# Remove those two rows entirely (shifting the remaining language
# rows up), which also shrinks the sheet from 111 to 109 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = $lastRow; $r -ge 1; $r--) {
    $val = $ws.Cells.Item($r, 1).Text
    if ($val -eq "history" -or $val -eq "check") {
        $ws.Rows.Item($r).EntireRow.Delete()
    }
}

# Match the saved view/selection state: scrolled down near the bottom
# of the now-shorter list, with B105 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 100
$win.ScrollColumn = 1
$ws.Range("B105").Select()
